$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.405.60'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '3.074.55'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.03'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.62'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.32%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").Value = '3.071.69'
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.80'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.28'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000240'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.23%  '
$ws.Range("D15").Value = '3.584.61'
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("D18").Value = '63.383.76'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = '3.076.57'
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.34'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  +1.19%  '
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.48'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.85'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.96'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.31'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").Value = '  -2.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.114'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.10'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.04'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("E39").Value = '  -4.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.60'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.21'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '440.96'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.21%  '
$ws.Range("E43").Value = '  -4.06%  '
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.82'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("D47").Value = '2.792.52'
$ws.Range("E47").Value = '  -4.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.59'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.51%  '
$ws.Range("E51").Value = '  +0.10%  '
